$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1895.5
$ws.Range("J70").Value = 1867.5
$ws.Range("L70").Value = 5602.5
$ws.Range("N70").Value = -6142.5
$ws.Range("H73").Value = 1895.5
$ws.Range("J73").Value = 1867.5
$ws.Range("L73").Value = 5602.5
$ws.Range("N73").Value = -7474.5
$ws.Range("H135").Value = 554.125
$ws.Range("I135").Value = 554.125
$ws.Range("K135").Value = 4987.125
$ws.Range("M135").Value = -2452.125
$ws.Range("H137").Value = 1897.5358
$ws.Range("I137").Value = 1553.8667
$ws.Range("K137").Value = 4661.6001
$ws.Range("M137").Value = -2111.6001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 333.18182
$ws.Range("I5").Value = 352.77777
$ws.Range("J5").Value = 245
$ws.Range("K5").Value = 352.77777
$ws.Range("L5").Value = 245
$ws.Range("M5").Value = -240.77777
$ws.Range("N5").Value = -469
$ws.Range("H19").Value = 2786.75
$ws.Range("I19").Value = 2049
$ws.Range("K19").Value = 2049
$ws.Range("M19").Value = -1820
$ws.Range("H32").Value = 3840.1892
$ws.Range("I32").Value = 3840.1892
$ws.Range("K32").Value = 3840.1892
$ws.Range("M32").Value = -3553.1892
$ws.Range("H39").Value = 14333
$ws.Range("I39").Value = 14333
$ws.Range("K39").Value = 14333
$ws.Range("M39").Value = -13813
$ws.Range("H61").Value = 996.3333
$ws.Range("J61").Value = 995
$ws.Range("L61").Value = 995
$ws.Range("N61").Value = -1419
$ws.Range("H74").Value = 1621.8
$ws.Range("I74").Value = 1703.0667
$ws.Range("J74").Value = 1499.9
$ws.Range("K74").Value = 1703.0667
$ws.Range("L74").Value = 1499.9
$ws.Range("M74").Value = -829.0667000000001
$ws.Range("N74").Value = -3247.9
$ws.Range("H77").Value = 1621.8
$ws.Range("I77").Value = 1703.0667
$ws.Range("J77").Value = 1499.9
$ws.Range("K77").Value = 8515.333500000001
$ws.Range("L77").Value = 7499.5
$ws.Range("M77").Value = -4147.333500000001
$ws.Range("N77").Value = -16235.5
$ws.Range("H88").Value = 2590
$ws.Range("J88").Value = 2521.8333
$ws.Range("L88").Value = 2521.8333
$ws.Range("N88").Value = -3333.8333
$ws.Range("H91").Value = 2590
$ws.Range("J91").Value = 2521.8333
$ws.Range("L91").Value = 2521.8333
$ws.Range("N91").Value = -5329.8333
$ws.Range("H132").Value = 1110.8695
$ws.Range("I132").Value = 1007.64703
$ws.Range("J132").Value = 1403.3334
$ws.Range("K132").Value = 3022.94109
$ws.Range("L132").Value = 4210.0002
$ws.Range("M132").Value = -492.9410899999998
$ws.Range("N132").Value = -9270.0002
$ws.Range("H136").Value = 996.3333
$ws.Range("J136").Value = 995
$ws.Range("L136").Value = 2985
$ws.Range("N136").Value = -8085

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 333.18182
$ws.Range("I4").Value = 352.77777
$ws.Range("J4").Value = 245
$ws.Range("K4").Value = 352.77777
$ws.Range("L4").Value = 245
$ws.Range("M4").Value = -237.77777
$ws.Range("N4").Value = -475
$ws.Range("H94").Value = 2509.4546
$ws.Range("I94").Value = 3159
$ws.Range("J94").Value = 1968.1666
$ws.Range("K94").Value = 3159
$ws.Range("L94").Value = 1968.1666
$ws.Range("M94").Value = -2708
$ws.Range("N94").Value = -2870.1666
$ws.Range("H134").Value = 3364.9
$ws.Range("I134").Value = 3364.9
$ws.Range("K134").Value = 10094.7
$ws.Range("M134").Value = -7559.700000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 21792
$ws.Range("J18").Value = 21792
$ws.Range("L18").Value = 21792
$ws.Range("N18").Value = -22252
$ws.Range("H31").Value = 2318.182
$ws.Range("J31").Value = 2128.0625
$ws.Range("L31").Value = 2128.0625
$ws.Range("N31").Value = -2718.0625
$ws.Range("H34").Value = 2318.182
$ws.Range("J34").Value = 2128.0625
$ws.Range("L34").Value = 2128.0625
$ws.Range("N34").Value = -2532.0625
$ws.Range("H122").Value = 415.5
$ws.Range("I122").Value = 415.5
$ws.Range("K122").Value = 1246.5
$ws.Range("M122").Value = 1203.5
$ws.Range("H134").Value = 2349.6538
$ws.Range("I134").Value = 2233
$ws.Range("J134").Value = 3749.5
$ws.Range("K134").Value = 6699
$ws.Range("L134").Value = 11248.5
$ws.Range("M134").Value = -4164
$ws.Range("N134").Value = -16318.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 267.66666
$ws.Range("I5").Value = 267.66666
$ws.Range("K5").Value = 802.9999799999999
$ws.Range("M5").Value = -690.9999799999999
$ws.Range("H135").Value = 267.66666
$ws.Range("I135").Value = 267.66666
$ws.Range("K135").Value = 2408.99994
$ws.Range("M135").Value = 126.0000600000003

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1266.1538
$ws.Range("I102").Value = 1097.1
$ws.Range("K102").Value = 1097.1
$ws.Range("M102").Value = 524.9000000000001
$ws.Range("H126").Value = 6869.9443
$ws.Range("I126").Value = 3221.5833
$ws.Range("K126").Value = 9664.749899999999
$ws.Range("M126").Value = -7194.749899999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2685.45
$ws.Range("I46").Value = 793.36365
$ws.Range("J46").Value = 3403.138
$ws.Range("K46").Value = 793.36365
$ws.Range("L46").Value = 3403.138
$ws.Range("M46").Value = -605.36365
$ws.Range("N46").Value = -3779.138
$ws.Range("H132").Value = 6628.6924
$ws.Range("I132").Value = 6141.4287
$ws.Range("K132").Value = 18424.2861
$ws.Range("M132").Value = -15894.2861
$ws.Range("H136").Value = 43480660
$ws.Range("I136").Value = 2979.0625
$ws.Range("J136").Value = 142858200
$ws.Range("K136").Value = 8937.1875
$ws.Range("L136").Value = 428574600
$ws.Range("M136").Value = -6387.1875
$ws.Range("N136").Value = -428579700

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1634.5
$ws.Range("I6").Value = 1002.3333
$ws.Range("J6").Value = 2266.6667
$ws.Range("K6").Value = 1002.3333
$ws.Range("L6").Value = 2266.6667
$ws.Range("M6").Value = -887.3333
$ws.Range("N6").Value = -2496.6667
$ws.Range("H7").Value = 2669.6667
$ws.Range("I7").Value = 2669.6667
$ws.Range("K7").Value = 2669.6667
$ws.Range("M7").Value = -2556.6667
$ws.Range("H9").Value = 35999
$ws.Range("I9").Value = 35999
$ws.Range("K9").Value = 35999
$ws.Range("M9").Value = -35859
$ws.Range("H30").Value = 24381.6
$ws.Range("J30").Value = 24977
$ws.Range("L30").Value = 24977
$ws.Range("N30").Value = -25191
$ws.Range("H113").Value = 609.0833
$ws.Range("I113").Value = 543.7273
$ws.Range("J113").Value = 1328
$ws.Range("K113").Value = 1631.1819
$ws.Range("L113").Value = 3984
$ws.Range("M113").Value = 538.8181
$ws.Range("N113").Value = -8324
$ws.Range("H122").Value = 3870.261
$ws.Range("I122").Value = 3637.0908
$ws.Range("K122").Value = 10911.2724
$ws.Range("M122").Value = -8461.2724
$ws.Range("H132").Value = 3529.3333
$ws.Range("I132").Value = 2383.3333
$ws.Range("J132").Value = 5248.3335
$ws.Range("K132").Value = 7149.999899999999
$ws.Range("L132").Value = 15745.0005
$ws.Range("M132").Value = -4619.999899999999
$ws.Range("N132").Value = -20805.0005
$ws.Range("H136").Value = 985.1070999999999
$ws.Range("I136").Value = 929
$ws.Range("K136").Value = 2787
$ws.Range("M136").Value = -237
